$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NATMI LR-pairs (App -> Cd74): Ligand/Receptor-expressing cell counts
# increase from 1 to 3 (cols E, K) for every data row (2-17); all
# downstream-derived columns (G,H,I,J,M,N,O,P,Q,R,S,T) are updated to
# the recomputed values that follow from that count change.

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 103.4275383333333
$ws.Range("H2").Value = 310.282615
$ws.Range("I2").Value = 0.2485530285127421
$ws.Range("J2").Value = 0.2485530285127421
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.656454333333333
$ws.Range("N2").Value = 4.969363
$ws.Range("O2").Value = 0.0006990759915034363
$ws.Range("P2").Value = 0.0006990759915034364
$ws.Range("Q2").Value = 171.3229940582495
$ws.Range("R2").Value = 1541.906946524245
$ws.Range("S2").Value = 0.0001737574548487271
$ws.Range("T2").Value = 0.0001737574548487271

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 103.4275383333333
$ws.Range("H3").Value = 310.282615
$ws.Range("I3").Value = 0.2485530285127421
$ws.Range("J3").Value = 0.2485530285127421
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.599386
$ws.Range("N3").Value = 19.798158
$ws.Range("O3").Value = 0.002785149109411345
$ws.Range("P3").Value = 0.002785149109411345
$ws.Range("Q3").Value = 682.5582484914634
$ws.Range("R3").Value = 6143.024236423171
$ws.Range("S3").Value = 0.0006922572460037563
$ws.Range("T3").Value = 0.0006922572460037563

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 103.4275383333333
$ws.Range("H4").Value = 310.282615
$ws.Range("I4").Value = 0.2485530285127421
$ws.Range("J4").Value = 0.2485530285127421
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2347.576497333333
$ws.Range("N4").Value = 7042.729492
$ws.Range("O4").Value = 0.9907513503260664
$ws.Range("P4").Value = 0.9907513503260664
$ws.Range("Q4").Value = 242804.0581683758
$ws.Range("R4").Value = 2185236.523515382
$ws.Range("S4").Value = 0.2462542486266325
$ws.Range("T4").Value = 0.2462542486266325

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 103.4275383333333
$ws.Range("H5").Value = 310.282615
$ws.Range("I5").Value = 0.2485530285127421
$ws.Range("J5").Value = 0.2485530285127421
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 13.65875266666666
$ws.Range("N5").Value = 40.97625799999999
$ws.Range("O5").Value = 0.005764424573018837
$ws.Range("P5").Value = 0.005764424573018838
$ws.Range("Q5").Value = 1412.691165017186
$ws.Range("R5").Value = 12714.22048515467
$ws.Range("S5").Value = 0.001432765185257102
$ws.Range("T5").Value = 0.001432765185257102

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 216.130539
$ws.Range("H6").Value = 648.391617
$ws.Range("I6").Value = 0.5193964865470273
$ws.Range("J6").Value = 0.5193964865470272
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.656454333333333
$ws.Range("N6").Value = 4.969363
$ws.Range("O6").Value = 0.0006990759915034363
$ws.Range("P6").Value = 0.0006990759915034364
$ws.Range("Q6").Value = 358.0103678922191
$ws.Range("R6").Value = 3222.093311029971
$ws.Range("S6").Value = 0.0003630976138162643
$ws.Range("T6").Value = 0.0003630976138162643

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 216.130539
$ws.Range("H7").Value = 648.391617
$ws.Range("I7").Value = 0.5193964865470273
$ws.Range("J7").Value = 0.5193964865470272
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.599386
$ws.Range("N7").Value = 19.798158
$ws.Range("O7").Value = 0.002785149109411345
$ws.Range("P7").Value = 0.002785149109411345
$ws.Range("Q7").Value = 1426.328853249054
$ws.Range("R7").Value = 12836.95967924149
$ws.Range("S7").Value = 0.001446596661937835
$ws.Range("T7").Value = 0.001446596661937834

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 216.130539
$ws.Range("H8").Value = 648.391617
$ws.Range("I8").Value = 0.5193964865470273
$ws.Range("J8").Value = 0.5193964865470272
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2347.576497333333
$ws.Range("N8").Value = 7042.729492
$ws.Range("O8").Value = 0.9907513503260664
$ws.Range("P8").Value = 0.9907513503260664
$ws.Range("Q8").Value = 507382.9737123854
$ws.Range("R8").Value = 4566446.763411469
$ws.Range("S8").Value = 0.5145927704010819
$ws.Range("T8").Value = 0.5145927704010818

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 216.130539
$ws.Range("H9").Value = 648.391617
$ws.Range("I9").Value = 0.5193964865470273
$ws.Range("J9").Value = 0.5193964865470272
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.65875266666666
$ws.Range("N9").Value = 40.97625799999999
$ws.Range("O9").Value = 0.005764424573018837
$ws.Range("P9").Value = 0.005764424573018838
$ws.Range("Q9").Value = 2952.073575914354
$ws.Range("R9").Value = 26568.66218322918
$ws.Range("S9").Value = 0.002994021870191332
$ws.Range("T9").Value = 0.002994021870191332

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 71.607325
$ws.Range("H10").Value = 214.821975
$ws.Range("I10").Value = 0.1720839321833696
$ws.Range("J10").Value = 0.1720839321833696
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.656454333333333
$ws.Range("N10").Value = 4.969363
$ws.Range("O10").Value = 0.0006990759915034363
$ws.Range("P10").Value = 0.0006990759915034364
$ws.Range("Q10").Value = 118.6142637946583
$ws.Range("R10").Value = 1067.528374151925
$ws.Range("S10").Value = 0.0001202997455128992
$ws.Range("T10").Value = 0.0001202997455128992

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 71.607325
$ws.Range("H11").Value = 214.821975
$ws.Range("I11").Value = 0.1720839321833696
$ws.Range("J11").Value = 0.1720839321833696
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 6.599386
$ws.Range("N11").Value = 19.798158
$ws.Range("O11").Value = 0.002785149109411345
$ws.Range("P11").Value = 0.002785149109411345
$ws.Range("Q11").Value = 472.56437810245
$ws.Range("R11").Value = 4253.079402922051
$ws.Range("S11").Value = 0.000479279410464514
$ws.Range("T11").Value = 0.0004792794104645141

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 71.607325
$ws.Range("H12").Value = 214.821975
$ws.Range("I12").Value = 0.1720839321833696
$ws.Range("J12").Value = 0.1720839321833696
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 2347.576497333333
$ws.Range("N12").Value = 7042.729492
$ws.Range("O12").Value = 0.9907513503260664
$ws.Range("P12").Value = 0.9907513503260664
$ws.Range("Q12").Value = 168103.6732069096
$ws.Range("R12").Value = 1512933.058862187
$ws.Range("S12").Value = 0.1704923881800926
$ws.Range("T12").Value = 0.1704923881800926

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 71.607325
$ws.Range("H13").Value = 214.821975
$ws.Range("I13").Value = 0.1720839321833696
$ws.Range("J13").Value = 0.1720839321833696
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 13.65875266666666
$ws.Range("N13").Value = 40.97625799999999
$ws.Range("O13").Value = 0.005764424573018837
$ws.Range("P13").Value = 0.005764424573018838
$ws.Range("Q13").Value = 978.0667412966166
$ws.Range("R13").Value = 8802.600671669548
$ws.Range("S13").Value = 0.0009919648472995226
$ws.Range("T13").Value = 0.0009919648472995228

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 24.953198
$ws.Range("H14").Value = 74.859594
$ws.Range("I14").Value = 0.05996655275686102
$ws.Range("J14").Value = 0.05996655275686102
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.656454333333333
$ws.Range("N14").Value = 4.969363
$ws.Range("O14").Value = 0.0006990759915034363
$ws.Range("P14").Value = 0.0006990759915034364
$ws.Range("Q14").Value = 41.33383295762467
$ws.Range("R14").Value = 372.0044966186221
$ws.Range("S14").Value = 0.00004192117732554574
$ws.Range("T14").Value = 0.00004192117732554574

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 24.953198
$ws.Range("H15").Value = 74.859594
$ws.Range("I15").Value = 0.05996655275686102
$ws.Range("J15").Value = 0.05996655275686102
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 6.599386
$ws.Range("N15").Value = 19.798158
$ws.Range("O15").Value = 0.002785149109411345
$ws.Range("P15").Value = 0.002785149109411345
$ws.Range("Q15").Value = 164.675785536428
$ws.Range("R15").Value = 1482.082069827852
$ws.Range("S15").Value = 0.0001670157910052399
$ws.Range("T15").Value = 0.0001670157910052399

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 24.953198
$ws.Range("H16").Value = 74.859594
$ws.Range("I16").Value = 0.05996655275686102
$ws.Range("J16").Value = 0.05996655275686102
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 2347.576497333333
$ws.Range("N16").Value = 7042.729492
$ws.Range("O16").Value = 0.9907513503260664
$ws.Range("P16").Value = 0.9907513503260664
$ws.Range("Q16").Value = 58579.54115810514
$ws.Range("R16").Value = 527215.8704229463
$ws.Range("S16").Value = 0.05941194311825936
$ws.Range("T16").Value = 0.05941194311825936

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 24.953198
$ws.Range("H17").Value = 74.859594
$ws.Range("I17").Value = 0.05996655275686102
$ws.Range("J17").Value = 0.05996655275686102
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 13.65875266666666
$ws.Range("N17").Value = 40.97625799999999
$ws.Range("O17").Value = 0.005764424573018837
$ws.Range("P17").Value = 0.005764424573018838
$ws.Range("Q17").Value = 340.8295597243613
$ws.Range("R17").Value = 3067.466037519252
$ws.Range("S17").Value = 0.0003456726702708802
$ws.Range("T17").Value = 0.0003456726702708802
